$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: append a new sentence at the end of the introductory paragraph,
# right after "...informationsmängderna analyseras utifrån era
# förutsättningar." Typing at the very end of the run creates a brand new
# run for the freshly-typed text (matches the diff: a new <w:r> with
# xml:space="preserve").
# ---------------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("informationsmängderna analyseras utifrån era förutsättningar.", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $r1.Collapse(0)
    $r1.InsertAfter(" Arkitekturen för kontraktets användning är den samma som för övriga kontrakt i denna tjänstedomän. Den finns beskriven i inledande stycken i tjänstekontraktsbeskrivningen (dock inte i detta utsnitt).")
}

# ---------------------------------------------------------------------------
# Edit 2: rename the "Email" table header to "Epost". Only the very first
# occurrence in the document is the standalone word we want; restrict the
# search to a narrow window around it so "orgUnitEmail" elsewhere is left
# untouched.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Start = 0
$found2 = $r2.Find.Execute("Email", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $r2.Text = "Epost"
}

# ---------------------------------------------------------------------------
# Edit 3: "Generella kommentarer rörande informationsmängden" becomes
# "Generella kommentarer rörande tjänstekontraktet".
# ---------------------------------------------------------------------------
$r3 = $d.Content
$found3 = $r3.Find.Execute("Generella kommentarer rörande informationsmängden", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Generella kommentarer rörande tjänstekontraktet", 2)

# ---------------------------------------------------------------------------
# Edit 4: "... returnerar ställda diagnoser för patienter, lagrade i
# journalsystem." becomes "... returnerar journalförda ställda diagnoser
# för patienter."
# ---------------------------------------------------------------------------
$r4 = $d.Content
$found4 = $r4.Find.Execute(" returnerar ställda diagnoser för patienter, lagrade i journalsystem.", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    " returnerar journalförda ställda diagnoser för patienter.", 2)

# ---------------------------------------------------------------------------
# Edit 5: "(ex. SOAP-header)." becomes "(ex. SOAP-header)." with the word
# order of "SOAP" swapped across the hyphen: " (ex. SOAP-" / "header" / ")."
# turns into " (ex. " / "SOAP" / "-header)." -- edited in-place, scoped to
# the small window right after the unique word "kuvertering" so the very
# common ")." token elsewhere in the document is not touched.
# ---------------------------------------------------------------------------
$anchor5 = $d.Content
$foundAnchor5 = $anchor5.Find.Execute("kuvertering", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundAnchor5) {
    $anchor5.Collapse(0)
    $winStart5 = $anchor5.Start

    $w1 = $d.Range($winStart5, $winStart5 + 40)
    $f1 = $w1.Find.Execute(" (ex. SOAP-", $false, $false, $false, $false, $false, $true, 1, $false, " (ex. ", 1)

    $w2 = $d.Range($winStart5, $winStart5 + 40)
    $f2 = $w2.Find.Execute("header", $false, $false, $false, $false, $false, $true, 1, $false, "SOAP", 1)

    $w3 = $d.Range($winStart5, $winStart5 + 40)
    $f3 = $w3.Find.Execute(").", $false, $false, $false, $false, $false, $true, 1, $false, "-header).", 1)
}
